# Regenerate orders with updated distance/size codes.
# Mapping applied to text values (substring replacement, order matters
# so that e.g. "D80" is not mistaken for anything else):
#   D80 -> D86
#   D51 -> D55
#   D64 -> D69
#   S30 -> S31
# (S20 / S25 remain unchanged)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
$lastCol = $usedRange.Columns.Count

# Columns whose text values may contain the distance/size codes:
#   B = Condition, D = Filename_Left, E = Filename_Right,
#   H = Distance, J = Size
$targetCols = @(2, 4, 5, 8, 10)

for ($r = 1; $r -le $lastRow; $r++) {
    foreach ($c in $targetCols) {
        $cell = $ws.Cells.Item($r, $c)
        # NOTE: bare `.Value` (no parens) does not reliably return the
        # scalar in this runtime - use `.Value2` (or `.Value()`) to read.
        $val = $cell.Value2
        if ($val -ne $null -and $val -is [string]) {
            $newVal = $val
            $newVal = $newVal.Replace("D80", "D86")
            $newVal = $newVal.Replace("D51", "D55")
            $newVal = $newVal.Replace("D64", "D69")
            $newVal = $newVal.Replace("S30", "S31")
            if ($newVal -ne $val) {
                $cell.Value = $newVal
            }
        }
    }
}
